# Append two new match rows (118 and 119) to the Romania Liga-1 2023-2024 sheet,
# matching the style/format already used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (117) down onto
# the two new rows so that column A (bold/bordered index) and column E
# (date-time number format) keep the same look as the rest of the table.
$ws.Range("A117:V117").Copy()
$ws.Range("A118:V119").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---- Row 118 : Farul Constanta 1 x 1 FC Botosani ----
$r = 118
$ws.Cells.Item($r, 1).Value = 117
$ws.Cells.Item($r, 2).Value = "romania"
$ws.Cells.Item($r, 3).Value = "liga-1"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45235.6875
$ws.Cells.Item($r, 6).Value = "Farul Constanta"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = "FC Botosani"
$ws.Cells.Item($r, 9).Value = 1

$ws.Cells.Item($r, 10).Value = 1.5
$ws.Cells.Item($r, 11).Value = "02/11/2023 14:42"
$ws.Cells.Item($r, 12).Value = 1.48
$ws.Cells.Item($r, 13).Value = "05/11/2023 16:20"

$ws.Cells.Item($r, 14).Value = 4.21
$ws.Cells.Item($r, 15).Value = "02/11/2023 14:42"
$ws.Cells.Item($r, 16).Value = 4.46
$ws.Cells.Item($r, 17).Value = "05/11/2023 16:29"

$ws.Cells.Item($r, 18).Value = 6.46
$ws.Cells.Item($r, 19).Value = "02/11/2023 14:42"
$ws.Cells.Item($r, 20).Value = 6.74
$ws.Cells.Item($r, 21).Value = "05/11/2023 16:29"

$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/farul-constanta-fc-botosani/Y5tXgOAI/"

# ---- Row 119 : FCSB 1 x 2 FC Rapid Bucuresti ----
$r = 119
$ws.Cells.Item($r, 1).Value = 118
$ws.Cells.Item($r, 2).Value = "romania"
$ws.Cells.Item($r, 3).Value = "liga-1"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45235.8125
$ws.Cells.Item($r, 6).Value = "FCSB"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = "FC Rapid Bucuresti"
$ws.Cells.Item($r, 9).Value = 2

$ws.Cells.Item($r, 10).Value = 1.99
$ws.Cells.Item($r, 11).Value = "02/11/2023 14:42"
$ws.Cells.Item($r, 12).Value = 1.85
$ws.Cells.Item($r, 13).Value = "05/11/2023 19:08"

$ws.Cells.Item($r, 14).Value = 3.51
$ws.Cells.Item($r, 15).Value = "02/11/2023 14:42"
$ws.Cells.Item($r, 16).Value = 3.75
$ws.Cells.Item($r, 17).Value = "05/11/2023 19:30"

$ws.Cells.Item($r, 18).Value = 3.75
$ws.Cells.Item($r, 19).Value = "02/11/2023 14:42"
$ws.Cells.Item($r, 20).Value = 4.14
$ws.Cells.Item($r, 21).Value = "05/11/2023 19:24"

$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/fcsb-rapid-bucuresti/bZYKdQQa/"
